# Add bot interface and section 'Aprender'
# - Clears the old sample data rows (2-7), keeping just the header row.
# - Adds three new header columns: "number", "ultima mensagem",
#   "Frases  / Aprender".
# - Moves the active selection to A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old sample rows (2-7), keep header row 1 intact.
# -4162 = xlShiftUp
$ws.Range("A2:D7").Delete(-4162)

# New header cells for the bot interface / "Aprender" section.
$ws.Range("E1").Value = "number"
$ws.Range("F1").Value = "ultima mensagem"
$ws.Range("G1").Value = "Frases  / Aprender"

# Give the new header cells the same (centered) formatting as the rest of
# row 1, without inventing a brand-new style entry.
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths for the newly introduced columns.
$ws.Columns.Item(5).ColumnWidth = 15.5703125
$ws.Columns.Item(6).ColumnWidth = 61.5703125
$ws.Columns.Item(7).ColumnWidth = 18.5703125

# Move the active cell/selection as recorded after the edit.
$ws.Range("A9").Select()
